# fix: change typo in functional requirements
#
# Rewrites the "Functional Requirements" bullet list so each item states a
# full requirement sentence instead of a short label, and removes the
# now-redundant trailing "Show analysis and transcript" bullet (its content
# was folded into the "Payment" bullet).

$d = $word.ActiveDocument

$APOS = [char]0x2019   # RIGHT SINGLE QUOTATION MARK (user's / user's)

function Split-RunAt {
    param($SearchRange, $Needle)

    # Locate $Needle inside $SearchRange and force the host to materialize
    # it as its own run (without altering its visible formatting) by
    # toggling Bold on/off across just that sub-range.
    $hit = $SearchRange.Duplicate
    $hit.Find.Execute($Needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if ($hit.Find.Found) {
        $hit.Bold = 1
        $hit.Bold = 0
    }
}

function Set-ParagraphRuns {
    param($Paragraph, [string[]] $Parts)

    $full = [string]::Join("", $Parts)

    # Replace the paragraph's text (but not its end-of-paragraph mark) with
    # the full new sentence as a single run first.
    $r = $Paragraph.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $full

    # Now split that single run into several runs at the boundaries implied
    # by $Parts, by toggling Bold on/off across each sub-range in turn.
    # Toggling a formatting property and then reverting it forces the host
    # to materialize a separate run for that sub-range while leaving the
    # run's properties identical to its neighbours.
    $pStart = $Paragraph.Range.Start
    $offset = 0
    for ($i = 0; $i -lt $Parts.Count; $i++) {
        $len = $Parts[$i].Length
        if ($len -gt 0 -and $i -lt ($Parts.Count - 1)) {
            $sub = $d.Range($pStart + $offset, $pStart + $offset + $len)
            $sub.Bold = 1
            $sub.Bold = 0
        }
        $offset += $len
    }
}

# Split the two technology names out of their surrounding sentences into
# their own runs (the text itself does not change), mirroring how Word
# ring-fences words the spell checker flagged with proofErr markers.
$paraTech = $d.Paragraphs.Item(5)
Split-RunAt $paraTech.Range "VueJs"
Split-RunAt $paraTech.Range "RestAPI"

# 1) Register
Set-ParagraphRuns $d.Paragraphs.Item(7) @(
    "The system must allow users to register with their email and password."
)

# 2) E-mail verification
Set-ParagraphRuns $d.Paragraphs.Item(8) @(
    ("The system must send a verification link to user" + $APOS + "s "),
    " ",
    "email to verify it."
)

# 3) Login
Set-ParagraphRuns $d.Paragraphs.Item(9) @(
    "The system must allow users to log in with their credentials (email and password)"
)

# 4) Forget password
Set-ParagraphRuns $d.Paragraphs.Item(10) @(
    ("The system must sent a link to user" + $APOS + "s mail to renew user" + $APOS + "s password"),
    "."
)

# 5) Renew password
Set-ParagraphRuns $d.Paragraphs.Item(11) @(
    "Show a list of the patients",
    ": The system must show patient list to the psychologist"
)

# 6) Show a list of the patients
Set-ParagraphRuns $d.Paragraphs.Item(12) @(
    "Show a detailed information of a patient",
    ": The system mush show patients detailed data to the user."
)

# 7) Show a detailed information of a patient
Set-ParagraphRuns $d.Paragraphs.Item(13) @(
    "The system must allow to user to",
    " ",
    "update their ",
    "patients data"
)

# 8) Update a specific patient
Set-ParagraphRuns $d.Paragraphs.Item(14) @(
    "The system must allow user to add or remove patient"
)

# 9) Add or remove a patient
Set-ParagraphRuns $d.Paragraphs.Item(15) @(
    "The system must allow user to d",
    "ownload transcript of the session"
)

# 10) Download transcript of the session
Set-ParagraphRuns $d.Paragraphs.Item(16) @(
    "The system must allow user to",
    " ",
    "Download analysis of the session"
)

# 11) Download analysis of the session
Set-ParagraphRuns $d.Paragraphs.Item(17) @(
    "The system must allow user to",
    " ",
    "pay for the service"
)

# 12) Payment
Set-ParagraphRuns $d.Paragraphs.Item(18) @(
    "The system must allow user to",
    " ",
    "see ",
    "analysis and transcript"
)

# 13) Show analysis and transcript -> removed entirely (its text was merged
# into the "Payment" bullet above).
$d.Paragraphs.Item(19).Range.Delete() | Out-Null
